$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells below contain numeric-looking text (e.g. "44.038.15",
# "2.44") that must remain plain text, matching the workbook's existing
# inline-string price formatting. Force Text number format on just the D cells
# that are changing (skipping D20, which holds non-numeric subscript text) so
# Excel does not silently coerce them into numbers.
$priceTextRows = @(2, 3, 5, 6, 7, 9, 10, 11, 12, 14, 15, 16, 17, 18, 19, 21, 22, 23, 25, 26, 27, 28, 29, 30, 33, 34, 35, 36, 37, 38, 39, 42, 43, 44, 45, 47, 48, 49, 50, 51)
foreach ($r in $priceTextRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "44.038.15"
$ws.Cells.Item(2, 5).Value = "  +0.66%  "
$ws.Cells.Item(3, 4).Value = "2.265.70"
$ws.Cells.Item(3, 5).Value = "  -0.53%  "
$ws.Cells.Item(4, 5).Value = "  -0.24%  "
$ws.Cells.Item(5, 4).Value = "233.03"
$ws.Cells.Item(5, 5).Value = "  +0.70%  "
$ws.Cells.Item(6, 4).Value = "0.649"
$ws.Cells.Item(6, 5).Value = "  +3.29%  "
$ws.Cells.Item(7, 4).Value = "63.90"
$ws.Cells.Item(7, 5).Value = "  -1.19%  "
$ws.Cells.Item(8, 5).Value = "  -0.19%  "
$ws.Cells.Item(9, 4).Value = "0.451"
$ws.Cells.Item(9, 5).Value = "  +6.34%  "
$ws.Cells.Item(10, 4).Value = "0.0984"
$ws.Cells.Item(10, 5).Value = "  +1.37%  "
$ws.Cells.Item(11, 4).Value = "58.59"
$ws.Cells.Item(11, 5).Value = "  +1.07%  "
$ws.Cells.Item(12, 4).Value = "26.60"
$ws.Cells.Item(12, 5).Value = "  +1.17%  "
$ws.Cells.Item(13, 5).Value = "  +1.65%  "
$ws.Cells.Item(14, 4).Value = "2.600.58"
$ws.Cells.Item(14, 5).Value = "  -0.75%  "
$ws.Cells.Item(15, 4).Value = "15.65"
$ws.Cells.Item(15, 5).Value = "  -0.69%  "
$ws.Cells.Item(16, 4).Value = "6.14"
$ws.Cells.Item(16, 5).Value = "  +3.68%  "
$ws.Cells.Item(17, 4).Value = "0.838"
$ws.Cells.Item(17, 5).Value = "  +2.44%  "
$ws.Cells.Item(18, 4).Value = "2.268.36"
$ws.Cells.Item(18, 5).Value = "  -0.67%  "
$ws.Cells.Item(19, 4).Value = "43.949.90"
$ws.Cells.Item(19, 5).Value = "  +0.67%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0990"
$ws.Cells.Item(20, 5).Value = "  +3.90%  "
$ws.Cells.Item(21, 4).Value = "73.77"
$ws.Cells.Item(21, 5).Value = "  +0.50%  "
$ws.Cells.Item(22, 4).Value = "6.11"
$ws.Cells.Item(22, 5).Value = "  -0.87%  "
$ws.Cells.Item(23, 4).Value = "250.84"
$ws.Cells.Item(23, 5).Value = "  +0.22%  "
$ws.Cells.Item(24, 5).Value = "  +0.03%  "
$ws.Cells.Item(25, 4).Value = "2.44"
$ws.Cells.Item(25, 5).Value = "  -3.94%  "
$ws.Cells.Item(26, 4).Value = "3.34"
$ws.Cells.Item(26, 5).Value = "  +21.53%  "
$ws.Cells.Item(27, 4).Value = "2.23"
$ws.Cells.Item(27, 5).Value = "  -4.61%  "
$ws.Cells.Item(28, 4).Value = "9.91"
$ws.Cells.Item(28, 5).Value = "  -0.71%  "
$ws.Cells.Item(29, 4).Value = "173.60"
$ws.Cells.Item(29, 5).Value = "  +0.90%  "
$ws.Cells.Item(30, 4).Value = "22.05"
$ws.Cells.Item(30, 5).Value = "  +7.24%  "
$ws.Cells.Item(31, 5).Value = "  +0.28%  "
$ws.Cells.Item(32, 5).Value = "  +0.67%  "
$ws.Cells.Item(33, 4).Value = "0.126"
$ws.Cells.Item(33, 5).Value = "  +3.20%  "
$ws.Cells.Item(34, 4).Value = "4.94"
$ws.Cells.Item(34, 5).Value = "  +4.43%  "
$ws.Cells.Item(35, 4).Value = "0.0687"
$ws.Cells.Item(35, 5).Value = "  -1.36%  "
$ws.Cells.Item(36, 4).Value = "4.95"
$ws.Cells.Item(36, 5).Value = "  -5.74%  "
$ws.Cells.Item(37, 4).Value = "3.69"
$ws.Cells.Item(37, 5).Value = "  -2.29%  "
$ws.Cells.Item(38, 4).Value = "6.51"
$ws.Cells.Item(38, 5).Value = "  -4.42%  "
$ws.Cells.Item(39, 4).Value = "2.30"
$ws.Cells.Item(39, 5).Value = "  -2.15%  "
$ws.Cells.Item(40, 5).Value = "  +3.34%  "
$ws.Cells.Item(42, 4).Value = "8.73"
$ws.Cells.Item(42, 5).Value = "  +2.84%  "
$ws.Cells.Item(43, 4).Value = "0.000224"
$ws.Cells.Item(43, 5).Value = "  -1.44%  "
$ws.Cells.Item(44, 4).Value = "17.33"
$ws.Cells.Item(44, 5).Value = "  +2.36%  "
$ws.Cells.Item(45, 4).Value = "98.55"
$ws.Cells.Item(45, 5).Value = "  +0.67%  "
$ws.Cells.Item(46, 5).Value = "  -2.00%  "
$ws.Cells.Item(47, 4).Value = "0.0953"
$ws.Cells.Item(47, 5).Value = "  -1.15%  "
$ws.Cells.Item(48, 4).Value = "2.37"
$ws.Cells.Item(48, 5).Value = "  +1.35%  "
$ws.Cells.Item(49, 4).Value = "4.36"
$ws.Cells.Item(49, 5).Value = "  -6.69%  "
$ws.Cells.Item(50, 4).Value = "1.449.45"
$ws.Cells.Item(50, 5).Value = "  -2.48%  "
$ws.Cells.Item(51, 4).Value = "9.95"
$ws.Cells.Item(51, 5).Value = "  -9.61%  "
